$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 97, pushing the old row 97 (Fecha 2021-05-19 / 44335)
# down to row 98. The new row 97 will receive a copy of the (pre-edit) row 96 data,
# and row 96 itself is then refreshed with the new weekly price point.
$ws.Rows.Item(97).Insert()

# New row 97 = copy of the original row 96 values (before this week's update)
$ws.Cells.Item(97, 1).Value = 8
$ws.Cells.Item(97, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(97, 3).Value = "Coquimbo"
$ws.Cells.Item(97, 4).Value = 44418
$ws.Cells.Item(97, 5).Value = 4
$ws.Cells.Item(97, 6).Value = 100114007
$ws.Cells.Item(97, 7).Value = "Jengibre"
$ws.Cells.Item(97, 8).Value = "Sin especificar"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 500
$ws.Cells.Item(97, 11).Value = 14000
$ws.Cells.Item(97, 12).Value = 15000
$ws.Cells.Item(97, 13).Value = 14500
$ws.Cells.Item(97, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(97, 15).Value = "Perú"
$ws.Cells.Item(97, 16).Value = 1115
$ws.Cells.Item(97, 17).Value = 13
$ws.Cells.Item(97, 18).Value = "Hortaliza"

# Row 97's date cell keeps the same date number format as the other Fecha cells
$ws.Cells.Item(97, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Refresh row 96 with this week's new data point (Volumen 460, Fecha 2023-01-13)
$ws.Cells.Item(96, 4).Value = 44939
$ws.Cells.Item(96, 10).Value = 460
